# "cambios en documento de observaciones y gráficas"
# Update the performance-measurement values on the "Datos de rendimiento"
# sheet. The two chart sheets ("Gráfico1"/"Gráfico2") plot these cells
# directly, so their displayed values follow once the source cells change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos de rendimiento")

# First table (column C, rows 3-8) - "Tiempo [ms]"
$ws.Range("C3").Value = 510.608
$ws.Range("C4").Value = 4042.2269999999999
$ws.Range("C6").Value = 1221.357

# Second table (column C, rows 14-19)
$ws.Range("C14").Value = 1046.019

# Move the active cell/selection to E7 (was E8)
$ws.Range("E7").Select()
